# Apply budget updates to Sheet1:
#  - E22 becomes a formula that sums the "Changements site" sub-total rows
#    (D19:D21), recalculating the total from 179 to 229.
#  - F23 (grand total) recalculates automatically from 1413 to 1463 since it
#    references E22 via SUM(E5,E13,E22).
#  - Update the active selection to F21 (last touched cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E22").Formula = "=SUM(D19:D21)"

$ws.Range("F21").Select() | Out-Null
